$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 463.85715
$ws.Range("I5").Value = 399.5
$ws.Range("J5").Value = 549.6667
$ws.Range("K5").Value = 399.5
$ws.Range("L5").Value = 549.6667
$ws.Range("M5").Value = -284.5
$ws.Range("N5").Value = -779.6667

$ws.Range("H69").Value = 11843.368
$ws.Range("J69").Value = 12376.875
$ws.Range("L69").Value = 37130.625
$ws.Range("N69").Value = -38878.625

$ws.Range("H72").Value = 11843.368
$ws.Range("J72").Value = 12376.875
$ws.Range("L72").Value = 111391.875
$ws.Range("N72").Value = -120127.875

$ws.Range("H125").Value = 3723.6428
$ws.Range("I125").Value = 2333.6667
$ws.Range("J125").Value = 4766.125
$ws.Range("K125").Value = 21003.0003
$ws.Range("L125").Value = 42895.125
$ws.Range("M125").Value = -18543.0003
$ws.Range("N125").Value = -47815.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1489.2
$ws.Range("I2").Value = 1272.2727
$ws.Range("K2").Value = 1272.2727
$ws.Range("M2").Value = -1159.2727

$ws.Range("H32").Value = 11074.571
$ws.Range("I32").Value = 10739.821
$ws.Range("K32").Value = 10739.821
$ws.Range("M32").Value = -10452.821

$ws.Range("H45").Value = 3200.1155
$ws.Range("I45").Value = 2198.625
$ws.Range("J45").Value = 3645.2222
$ws.Range("K45").Value = 2198.625
$ws.Range("L45").Value = 3645.2222
$ws.Range("M45").Value = -1821.625
$ws.Range("N45").Value = -4399.2222

$ws.Range("H61").Value = 1958.3846
$ws.Range("I61").Value = 1216
$ws.Range("K61").Value = 1216
$ws.Range("M61").Value = -1004

$ws.Range("H63").Value = 1218.5555
$ws.Range("I63").Value = 996
$ws.Range("K63").Value = 996
$ws.Range("M63").Value = -310

$ws.Range("H66").Value = 1218.5555
$ws.Range("I66").Value = 996
$ws.Range("K66").Value = 4980
$ws.Range("M66").Value = -1548

$ws.Range("H97").Value = 3389.5715
$ws.Range("I97").Value = 899.5
$ws.Range("J97").Value = 6709.6665
$ws.Range("K97").Value = 899.5
$ws.Range("L97").Value = 6709.6665
$ws.Range("M97").Value = -403.5
$ws.Range("N97").Value = -7701.6665

$ws.Range("H116").Value = 1489.2
$ws.Range("I116").Value = 1272.2727
$ws.Range("K116").Value = 1272.2727
$ws.Range("M116").Value = 1021.7273

$ws.Range("H132").Value = 4215.7334
$ws.Range("I132").Value = 4088.3572
$ws.Range("K132").Value = 12265.0716
$ws.Range("M132").Value = -9735.071599999999

$ws.Range("H136").Value = 1958.3846
$ws.Range("I136").Value = 1216
$ws.Range("K136").Value = 3648
$ws.Range("M136").Value = -1098

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1489.2
$ws.Range("I3").Value = 1272.2727
$ws.Range("K3").Value = 1272.2727
$ws.Range("M3").Value = -1158.2727

$ws.Range("H86").Value = 3088.0454
$ws.Range("I86").Value = 2964.9375
$ws.Range("K86").Value = 2964.9375
$ws.Range("M86").Value = -1841.9375

$ws.Range("H89").Value = 3088.0454
$ws.Range("I89").Value = 2964.9375
$ws.Range("K89").Value = 14824.6875
$ws.Range("M89").Value = -9208.6875

$ws.Range("H99").Value = 853.4400000000001
$ws.Range("I99").Value = 813.5
$ws.Range("K99").Value = 813.5
$ws.Range("M99").Value = 684.5

$ws.Range("H105").Value = 2896.4
$ws.Range("I105").Value = 1884.1
$ws.Range("K105").Value = 1884.1
$ws.Range("M105").Value = -137.0999999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5867.1
$ws.Range("I86").Value = 4238.857
$ws.Range("K86").Value = 4238.857
$ws.Range("M86").Value = -3115.857

$ws.Range("H89").Value = 5867.1
$ws.Range("I89").Value = 4238.857
$ws.Range("K89").Value = 21194.285
$ws.Range("M89").Value = -15578.285

$ws.Range("H99").Value = 2071.2856
$ws.Range("J99").Value = 2328.1428
$ws.Range("L99").Value = 2328.1428
$ws.Range("N99").Value = -5324.1428

$ws.Range("H105").Value = 1612.8334
$ws.Range("I105").Value = 1645.909
$ws.Range("J105").Value = 1249
$ws.Range("K105").Value = 1645.909
$ws.Range("L105").Value = 1249
$ws.Range("M105").Value = 101.0909999999999
$ws.Range("N105").Value = -4743

$ws.Range("H122").Value = 3062
$ws.Range("I122").Value = 2586.8823
$ws.Range("K122").Value = 7760.646900000001
$ws.Range("M122").Value = -5310.646900000001

$ws.Range("H126").Value = 2071.2856
$ws.Range("J126").Value = 2328.1428
$ws.Range("L126").Value = 6984.428400000001
$ws.Range("N126").Value = -11924.4284

$ws.Range("H132").Value = 2699.4211
$ws.Range("I132").Value = 2130.3215
$ws.Range("K132").Value = 6390.9645
$ws.Range("M132").Value = -3860.9645

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1050
$ws.Range("I109").Value = 1050
$ws.Range("K109").Value = 3150
$ws.Range("M109").Value = -2110

$ws.Range("H117").Value = 3103.4285
$ws.Range("I117").Value = 729
$ws.Range("J117").Value = 3499.1667
$ws.Range("K117").Value = 2187
$ws.Range("L117").Value = 10497.5001
$ws.Range("M117").Value = 1255
$ws.Range("N117").Value = -17381.5001

$ws.Range("H132").Value = 3957.5518
$ws.Range("I132").Value = 2227
$ws.Range("J132").Value = 4318.0835
$ws.Range("K132").Value = 20043
$ws.Range("L132").Value = 38862.7515
$ws.Range("M132").Value = -17513
$ws.Range("N132").Value = -43922.7515

$ws.Range("H137").Value = 49100.332
$ws.Range("I137").Value = 75957.71000000001
$ws.Range("J137").Value = 11500
$ws.Range("K137").Value = 227873.13
$ws.Range("L137").Value = 34500
$ws.Range("M137").Value = -222773.13
$ws.Range("N137").Value = -44700

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5212.4287
$ws.Range("I126").Value = 4768.5
$ws.Range("K126").Value = 14305.5
$ws.Range("M126").Value = -11835.5

$ws.Range("H132").Value = 3060.2727
$ws.Range("I132").Value = 2442.121
$ws.Range("K132").Value = 7326.363
$ws.Range("M132").Value = -4796.363

$ws.Range("H136").Value = 55486
$ws.Range("J136").Value = 55486
$ws.Range("L136").Value = 166458
$ws.Range("N136").Value = -171558

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3317.4119
$ws.Range("J136").Value = 3436.842
$ws.Range("L136").Value = 10310.526
$ws.Range("N136").Value = -15410.526

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2431.9312
$ws.Range("I132").Value = 1941.48
$ws.Range("J132").Value = 5497.25
$ws.Range("K132").Value = 5824.440000000001
$ws.Range("L132").Value = 16491.75
$ws.Range("M132").Value = -3294.440000000001
$ws.Range("N132").Value = -21551.75

$ws.Range("H136").Value = 5513.477
$ws.Range("I136").Value = 4566.6943
$ws.Range("K136").Value = 13700.0829
$ws.Range("M136").Value = -11150.0829
